$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save" — copy formatting (style) from the neighboring
# header cell G1 (bold, centered, bordered) before writing the value so the
# resulting cell carries the same style index used by the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell H2: numeric value 0 (no special style, like the other data cells)
$ws.Range("H2").Value = 0
